$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2..271).
# Update every value from 45204 to 45205 (2023-10-05 -> 2023-10-06).
$ws.Range("C2:C271").Value = 45205
